$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above current row 2 - existing rows 2-57 shift down to 3-58
$ws.Rows.Item(2).Insert()

# Copy the formatting from the (now) row 3 down into the newly inserted row 2
# so the new row matches the look of the surrounding data rows.
$ws.Range("A3:P3").Copy()
$ws.Range("A2:P2").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Populate the newly inserted row 2 with the new product data
$ws.Range("A2").Value = 7791708001378
$ws.Range("B2").Value = "Pan"
$ws.Range("C2").Value = "para"
$ws.Range("D2").Value = "panchos"
$ws.Range("E2").Value = "Veneziana"
$ws.Range("F2").Value = 6
$ws.Range("G2").Value = "und."
$ws.Range("H2").Value = "bolsa"
$ws.Range("I2").Value = "Panes"
$ws.Range("J2").Value = "Argentina"
$ws.Range("K2").Value = 12
$ws.Range("L2").Value = $false
$ws.Range("M2").Value = $true
$ws.Range("N2").Value = "C:\VentaSoft\Imágenes de artículos\7791708001378.png"
$ws.Range("O2").Value = $true
$ws.Range("P2").Value = $true
